$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Range("H9").Value = 5920.55
$ws.Range("I9").Value = 671.2941
$ws.Range("K9").Value = 671.2941
$ws.Range("M9").Value = -502.2941

# Row 86
$ws.Range("H86").Value = 1658.1666
$ws.Range("I86").Value = 1885.1428
$ws.Range("J86").Value = 1340.4
$ws.Range("K86").Value = 1885.1428
$ws.Range("L86").Value = 1340.4
$ws.Range("M86").Value = -762.1428000000001
$ws.Range("N86").Value = -3586.4

# Row 88
$ws.Range("H88").Value = 1166.037
$ws.Range("J88").Value = 1031.1578
$ws.Range("L88").Value = 1031.1578
$ws.Range("N88").Value = -1843.1578

# Row 89
$ws.Range("H89").Value = 1658.1666
$ws.Range("I89").Value = 1885.1428
$ws.Range("J89").Value = 1340.4
$ws.Range("K89").Value = 9425.714
$ws.Range("L89").Value = 6702
$ws.Range("M89").Value = -3809.714
$ws.Range("N89").Value = -17934

# Row 91
$ws.Range("H91").Value = 1166.037
$ws.Range("J91").Value = 1031.1578
$ws.Range("L91").Value = 1031.1578
$ws.Range("N91").Value = -3839.1578

# Row 98
$ws.Range("H98").Value = 2117.5715
$ws.Range("I98").Value = 2284.4092
$ws.Range("K98").Value = 2284.4092
$ws.Range("M98").Value = -786.4092000000001

# Row 122
$ws.Range("H122").Value = 2117.5715
$ws.Range("I122").Value = 2284.4092
$ws.Range("K122").Value = 6853.2276
$ws.Range("M122").Value = -4403.2276

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 70701
$ws.Range("I86").Value = 1478.125
$ws.Range("K86").Value = 1478.125
$ws.Range("M86").Value = -355.125

# Row 89
$ws.Range("H89").Value = 70701
$ws.Range("I89").Value = 1478.125
$ws.Range("K89").Value = 7390.625
$ws.Range("M89").Value = -1774.625

# Row 141
$ws.Range("H141").Value = 190000
$ws.Range("J141").Value = 190000
$ws.Range("L141").Value = 190000
$ws.Range("N141").Value = -200360

$ws = $wb.Worksheets.Item("CRP")
# Row 29
$ws.Range("H29").Value = 13991.875
$ws.Range("I29").Value = 5509
$ws.Range("J29").Value = 16819.5
$ws.Range("K29").Value = 5509
$ws.Range("L29").Value = 16819.5
$ws.Range("M29").Value = -5216
$ws.Range("N29").Value = -17405.5

# Row 58
$ws.Range("H58").Value = 8098.2256
$ws.Range("I58").Value = 4424.8696
$ws.Range("J58").Value = 18659.125
$ws.Range("K58").Value = 4424.8696
$ws.Range("L58").Value = 18659.125
$ws.Range("M58").Value = -4221.8696
$ws.Range("N58").Value = -19065.125

# Row 122
$ws.Range("H122").Value = 992.8570999999999
$ws.Range("I122").Value = 788.4
$ws.Range("J122").Value = 1504
$ws.Range("K122").Value = 2365.2
$ws.Range("L122").Value = 4512
$ws.Range("M122").Value = 84.80000000000018
$ws.Range("N122").Value = -9412

# Row 134
$ws.Range("H134").Value = 1258.2307
$ws.Range("I134").Value = 840.2174
$ws.Range("K134").Value = 2520.6522
$ws.Range("M134").Value = 14.34780000000001

# Row 136
$ws.Range("H136").Value = 8098.2256
$ws.Range("I136").Value = 4424.8696
$ws.Range("J136").Value = 18659.125
$ws.Range("K136").Value = 13274.6088
$ws.Range("L136").Value = 55977.375
$ws.Range("M136").Value = -10724.6088
$ws.Range("N136").Value = -61077.375

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 861.3333
$ws.Range("I5").Value = 799
$ws.Range("J5").Value = 923.6667
$ws.Range("K5").Value = 2397
$ws.Range("L5").Value = 2771.0001
$ws.Range("M5").Value = -2285
$ws.Range("N5").Value = -2995.0001

# Row 45
$ws.Range("H45").Value = 1266.5
$ws.Range("J45").Value = 1266.5
$ws.Range("L45").Value = 3799.5
$ws.Range("N45").Value = -4863.5

# Row 107
$ws.Range("H107").Value = 1502.7391
$ws.Range("I107").Value = 386.5
$ws.Range("J107").Value = 1896.7059
$ws.Range("K107").Value = 1159.5
$ws.Range("L107").Value = 5690.1177
$ws.Range("M107").Value = 760.5
$ws.Range("N107").Value = -9530.117699999999

# Row 113
$ws.Range("H113").Value = 3973.3845
$ws.Range("I113").Value = 1651
$ws.Range("J113").Value = 4395.636
$ws.Range("K113").Value = 4953
$ws.Range("L113").Value = 13186.908
$ws.Range("M113").Value = -2783
$ws.Range("N113").Value = -17526.908

# Row 129
$ws.Range("H129").Value = 4797.346
$ws.Range("J129").Value = 5445.5
$ws.Range("L129").Value = 16336.5
$ws.Range("N129").Value = -26336.5

# Row 135
$ws.Range("H135").Value = 861.3333
$ws.Range("I135").Value = 799
$ws.Range("J135").Value = 923.6667
$ws.Range("K135").Value = 7191
$ws.Range("L135").Value = 8313.0003
$ws.Range("M135").Value = -4656
$ws.Range("N135").Value = -13383.0003

# Row 137
$ws.Range("H137").Value = 5342.071
$ws.Range("J137").Value = 5140.439
$ws.Range("L137").Value = 15421.317
$ws.Range("N137").Value = -25621.317

# Row 139
$ws.Range("H139").Value = 4178.579
$ws.Range("I139").Value = 5075.727
$ws.Range("J139").Value = 2945
$ws.Range("K139").Value = 15227.181
$ws.Range("L139").Value = 8835
$ws.Range("M139").Value = -10087.181
$ws.Range("N139").Value = -19115

$ws = $wb.Worksheets.Item("GSM")
# Row 5
$ws.Range("H5").Value = 580002
$ws.Range("I5").Value = 580002
$ws.Range("K5").Value = 580002
$ws.Range("M5").Value = -579890

# Row 102
$ws.Range("H102").Value = 34592.027
$ws.Range("I102").Value = 42816.42
$ws.Range("K102").Value = 42816.42
$ws.Range("M102").Value = -41194.42

# Row 109
$ws.Range("H109").Value = 31666
$ws.Range("J109").Value = 31666
$ws.Range("L109").Value = 31666
$ws.Range("N109").Value = -33746

# Row 126
$ws.Range("H126").Value = 43918.684
$ws.Range("I126").Value = 55798.59
$ws.Range("K126").Value = 167395.77
$ws.Range("M126").Value = -164925.77

$ws = $wb.Worksheets.Item("LTW")
# Row 82
$ws.Range("H82").Value = 2713.5715
$ws.Range("I82").Value = 2832.5
$ws.Range("J82").Value = 2000
$ws.Range("K82").Value = 2832.5
$ws.Range("L82").Value = 2000
$ws.Range("M82").Value = -2471.5
$ws.Range("N82").Value = -2722

# Row 85
$ws.Range("H85").Value = 2713.5715
$ws.Range("I85").Value = 2832.5
$ws.Range("J85").Value = 2000
$ws.Range("K85").Value = 2832.5
$ws.Range("L85").Value = 2000
$ws.Range("M85").Value = -1584.5
$ws.Range("N85").Value = -4496

# Row 100
$ws.Range("H100").Value = 1205.4546
$ws.Range("I100").Value = 1285.625
$ws.Range("K100").Value = 1285.625
$ws.Range("M100").Value = -744.625

# Row 132
$ws.Range("H132").Value = 4057
$ws.Range("I132").Value = 3410.4443
$ws.Range("K132").Value = 10231.3329
$ws.Range("M132").Value = -7701.332900000001

# Row 133
$ws.Range("H133").Value = 99999
$ws.Range("J133").Value = 99999
$ws.Range("L133").Value = 99999
$ws.Range("N133").Value = -105059

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 6416.263
$ws.Range("I132").Value = 3064.7856
$ws.Range("K132").Value = 9194.356800000001
$ws.Range("M132").Value = -6664.356800000001

# Row 136
$ws.Range("H136").Value = 4462.8813
$ws.Range("I136").Value = 4452.915
$ws.Range("K136").Value = 13358.745
$ws.Range("M136").Value = -10808.745
